# Refresh the crypto price/volume snapshot to match the latest scrape.
# (Coin/Link swap on rows 44-45: MXToken now ranks above RocketPoolETH.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.347.58'
$ws.Range('E2').Value = '  +4.07%  '

$ws.Range('D3').Value = '1.732.43'
$ws.Range('E3').Value = '  +2.54%  '

$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = "'219.39"
$ws.Range('E5').Value = '  +1.30%  '

$ws.Range('D6').Value = "'0.522"
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('E7').Value = '  -0.17%  '

$ws.Range('D8').Value = "'24.05"
$ws.Range('E8').Value = '  +4.12%  '

$ws.Range('E9').Value = '  +2.05%  '

$ws.Range('D10').Value = "'0.0638"
$ws.Range('E10').Value = '  +1.46%  '

$ws.Range('E11').Value = '  +0.47%  '

$ws.Range('D12').Value = '1.976.65'
$ws.Range('E12').Value = '  +2.59%  '

$ws.Range('D13').Value = '1.730.27'
$ws.Range('E13').Value = '  +3.18%  '

$ws.Range('E14').Value = '  +1.32%  '

$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('E16').Value = '  +0.59%  '

$ws.Range('D17').Value = '28.336.30'
$ws.Range('E17').Value = '  +4.12%  '

$ws.Range('D18').Value = "'247.42"
$ws.Range('E18').Value = '  +3.92%  '

$ws.Range('E19').Value = '  +1.12%  '

$ws.Range('D20').Value = "'7.92"

$ws.Range('E21').Value = '  -0.17%  '

$ws.Range('E22').Value = '  +1.31%  '

$ws.Range('E24').Value = '  -0.88%  '

$ws.Range('D25').Value = "'149.57"
$ws.Range('E25').Value = '  +0.75%  '

$ws.Range('E26').Value = '  +2.43%  '

$ws.Range('E27').Value = '  +1.16%  '

$ws.Range('E28').Value = '  +0.50%  '

$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('E30').Value = '  +2.59%  '

$ws.Range('E31').Value = '  +2.28%  '

$ws.Range('E32').Value = '  +0.50%  '

$ws.Range('D33').Value = "'3.28"
$ws.Range('E33').Value = '  +0.45%  '

$ws.Range('D34').Value = '1.491.54'
$ws.Range('E34').Value = '  -5.82%  '

$ws.Range('E35').Value = '  -2.32%  '

$ws.Range('D36').Value = "'0.980"
$ws.Range('E36').Value = '  +2.24%  '

$ws.Range('D37').Value = "'0.603"
$ws.Range('E37').Value = '  -0.29%  '

$ws.Range('E38').Value = '  +0.49%  '

$ws.Range('E39').Value = '  +1.24%  '

$ws.Range('E40').Value = '  +0.21%  '

$ws.Range('D41').Value = "'70.26"
$ws.Range('E41').Value = '  +0.84%  '

$ws.Range('E42').Value = '  -0.23%  '

$ws.Range('D43').Value = "'5.67"
$ws.Range('E43').Value = '  -2.41%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = "'2.30"
$ws.Range('E44').Value = '  +1.48%  '

$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.880.10'
$ws.Range('E45').Value = '  +2.24%  '

$ws.Range('E46').Value = '  +1.23%  '

$ws.Range('E47').Value = '  +7.08%  '

$ws.Range('E48').Value = '  +3.98%  '

$ws.Range('D49').Value = "'90.47"
$ws.Range('E49').Value = '  -0.99%  '

$ws.Range('D50').Value = "'8.22"
$ws.Range('E50').Value = '  -0.73%  '

$ws.Range('E51').Value = '  -1.06%  '

# Strip the auto-applied text-quote-prefix formatting so styles stay identical.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
